$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New item: necro hood (row 18)
# Note: shared-string table order must match the canonical diff, so
# write the Item Description (B) before the Item ID (A).
$ws.Cells.Item(18, 2).Value = "itd_head_necrot3"
$ws.Cells.Item(18, 1).Value = "it_eq_head_necrot3"
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = "0, 200"
$ws.Cells.Item(18, 5).Value = 120
$ws.Cells.Item(18, 6).Value = 10
$ws.Cells.Item(18, 13).Value = 10
$ws.Cells.Item(18, 15).Value = 5
$ws.Cells.Item(18, 18).Value = 10
$ws.Cells.Item(18, 19).Value = 10
$ws.Cells.Item(18, 25).Value = "res/assets/equipment/head/spritesheet_head_necrot3.png"

# Fix selection (time-bug) - move active selection to Y19
$ws.Range("Y19").Select() | Out-Null
